# participants_template.xlsx — add new registration fields to the header row
# and reflow the columns to make room for them.
#
# Old layout (A:J): id_no, name, gender, age, category, institution,
#                    nationality, phone, address, email
# New layout (A:M): id_no, name, gender, age, category, institution,
#                    institution_ownership, education_level, subjects,
#                    nationality, phone, district, email
# ("address" is dropped; institution_ownership / subjects / district /
#  education_level are new.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id_no"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "gender"
$ws.Range("D1").Value = "age"
$ws.Range("E1").Value = "category"
$ws.Range("F1").Value = "institution"

# New columns are written in the same left-to-right order they first appear
# in the header row so the workbook's shared-string table picks up the same
# "first seen" ordering as the authored workbook.
$ws.Range("G1").Value = "institution_ownership"
$ws.Range("I1").Value = "subjects"
$ws.Range("L1").Value = "district"
$ws.Range("H1").Value = "education_level"

$ws.Range("J1").Value = "nationality"
$ws.Range("K1").Value = "phone"
$ws.Range("M1").Value = "email"

# Widen/resize the columns to fit the new headers (G:M).
$ws.Columns.Item(7).ColumnWidth = 20.830729166666668
$ws.Columns.Item(8).ColumnWidth = 20.830729166666668
$ws.Columns.Item(9).ColumnWidth = 18.166666666666668
$ws.Columns.Item(10).ColumnWidth = 20.166666666666668
$ws.Columns.Item(11).ColumnWidth = 15.330729166666666
$ws.Columns.Item(12).ColumnWidth = 14.666666666666666
$ws.Columns.Item(13).ColumnWidth = 13.166666666666666

# Leave the sheet selection where the author last left it.
$ws.Range("E17").Select() | Out-Null
